$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet: Applications -> PrivateApp ---
$ws.Name = "PrivateApp"

# --- Wipe existing content/formatting so we can rebuild cleanly ---
$ws.Cells.Clear()

# --- Header row (row 1) ---
$ws.Cells.Item(1,1).Value = "app_name"
$ws.Cells.Item(1,2).Value = "host"
$ws.Cells.Item(1,3).Value = "port"
$ws.Cells.Item(1,4).Value = "protocol"
$ws.Cells.Item(1,5).Value = "publisher_id"
$ws.Cells.Item(1,6).Value = "publisher_name"
$ws.Cells.Item(1,7).Value = "tags"
$ws.Cells.Item(1,8).Value = "use_publisher_dns"
$ws.Cells.Item(1,9).Value = "clientless_access"
$ws.Cells.Item(1,10).Value = "private_app_protocol"

# --- Data rows ---
# Row 2: testeAPI2
$ws.Cells.Item(2,1).Value = "testeAPI2"
$ws.Cells.Item(2,2).Value = "webserver.local"
$ws.Cells.Item(2,3).Value = 80
$ws.Cells.Item(2,4).Value = "tcp"
$ws.Cells.Item(2,5).Value = 6
$ws.Cells.Item(2,6).Value = "Azure-Publisher"
$ws.Cells.Item(2,7).Value = "webserver"
$ws.Cells.Item(2,8).Value = "true"
$ws.Cells.Item(2,9).Value = "false"
$ws.Cells.Item(2,10).Value = "false"

# Row 3: APIteste
$ws.Cells.Item(3,1).Value = "APIteste"
$ws.Cells.Item(3,2).Value = "server.local"
$ws.Cells.Item(3,3).Value = 80
$ws.Cells.Item(3,4).Value = "tcp"
$ws.Cells.Item(3,5).Value = 6
$ws.Cells.Item(3,6).Value = "Azure-Publisher"
$ws.Cells.Item(3,7).Value = "servers"
$ws.Cells.Item(3,8).Value = "false"
$ws.Cells.Item(3,9).Value = "true"
$ws.Cells.Item(3,10).Value = "http"

# Row 4: Xx
$ws.Cells.Item(4,1).Value = "Xx"
$ws.Cells.Item(4,2).Value = "vdi.local"
$ws.Cells.Item(4,4).Value = "tcp"
$ws.Cells.Item(4,5).Value = 6
$ws.Cells.Item(4,6).Value = "Azure-Publisher"
$ws.Cells.Item(4,7).Value = "servers"
$ws.Cells.Item(4,8).Value = "false"
$ws.Cells.Item(4,9).Value = "false"
$ws.Cells.Item(4,10).Value = "false"
# port is literal text "80,443" -> force text format before assigning so it is not parsed as a number
$ws.Cells.Item(4,3).NumberFormat = "@"
$ws.Cells.Item(4,3).Value = "80,443"

# Row 5: AD01
$ws.Cells.Item(5,1).Value = "AD01"
$ws.Cells.Item(5,2).Value = "192.168.201.1"
$ws.Cells.Item(5,3).Value = 80
$ws.Cells.Item(5,4).Value = "tcp"
$ws.Cells.Item(5,5).Value = 6
$ws.Cells.Item(5,6).Value = "Azure-Publisher"
$ws.Cells.Item(5,7).Value = "machine"
$ws.Cells.Item(5,8).Value = "false"
$ws.Cells.Item(5,9).Value = "false"
$ws.Cells.Item(5,10).Value = "false"

# Row 6: AD02
$ws.Cells.Item(6,1).Value = "AD02"
$ws.Cells.Item(6,2).Value = "192.168.3.3"
$ws.Cells.Item(6,3).Value = 443
$ws.Cells.Item(6,4).Value = "tcp"
$ws.Cells.Item(6,5).Value = 6
$ws.Cells.Item(6,6).Value = "Azure-Publisher"
$ws.Cells.Item(6,7).Value = "machine"
$ws.Cells.Item(6,8).Value = "false"
$ws.Cells.Item(6,9).Value = "false"
$ws.Cells.Item(6,10).Value = "false"

# --- Number formats ---
# Column C (port) displayed as text, but keep numeric rows as real numbers (format applied after value)
$ws.Range("C1:C6").NumberFormat = "@"
# Column G (tags) text format for data rows
$ws.Range("G2:G6").NumberFormat = "@"
# Column J (header + "http" row) text format
$ws.Range("J1").NumberFormat = "@"
$ws.Range("J3").NumberFormat = "@"
# Column H (use_publisher_dns) boolean-like display format for data rows
$ws.Range("H2:H6").NumberFormat = """VERDADEIRO"";""VERDADEIRO"";""FALSO"""

# --- Fonts ---
# Green JetBrains Mono font (matches existing accent font) on A2
$ws.Range("A2").Font.Name = $ws.Range("F2").Font.Name
$ws.Range("A2").Font.Size = $ws.Range("F2").Font.Size
$ws.Range("A2").Font.Color = $ws.Range("F2").Font.Color

# Black JetBrains Mono font on B2, H1, I1:I6, J2, J4:J6
$ws.Range("B2").Font.Name = $ws.Range("F2").Font.Name
$ws.Range("B2").Font.Size = $ws.Range("F2").Font.Size
$ws.Range("B2").Font.Color = 0

$ws.Range("H1").Font.Name = $ws.Range("F2").Font.Name
$ws.Range("H1").Font.Size = $ws.Range("F2").Font.Size
$ws.Range("H1").Font.Color = 0

$ws.Range("I1:I6").Font.Name = $ws.Range("F2").Font.Name
$ws.Range("I1:I6").Font.Size = $ws.Range("F2").Font.Size
$ws.Range("I1:I6").Font.Color = 0

$ws.Range("J2").Font.Name = $ws.Range("F2").Font.Name
$ws.Range("J2").Font.Size = $ws.Range("F2").Font.Size
$ws.Range("J2").Font.Color = 0

$ws.Range("J4:J6").Font.Name = $ws.Range("F2").Font.Name
$ws.Range("J4:J6").Font.Size = $ws.Range("F2").Font.Size
$ws.Range("J4:J6").Font.Color = 0

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 40.92
$ws.Columns.Item(8).ColumnWidth = 22.25
$ws.Columns.Item(9).ColumnWidth = 18.09
$ws.Columns.Item(10).ColumnWidth = 19.75

# --- View state ---
$excel.ActiveWindow.Zoom = 110
$ws.Range("G6").Select()
